$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.573
$ws.Range("B4").Value = 6.837999999999999
$ws.Range("C4").Value = -12.309

$ws.Range("B5").Value = 6.226

$ws.Range("A7").Value = -21.23

$ws.Range("B8").Value = 6.256

$ws.Range("C9").Value = -11.854

$ws.Range("A16").Value = -21.312
$ws.Range("B16").Value = 5.999000000000001

$ws.Range("C18").Value = -12.732
